# SQL Saturday 600 title slide update:
#   "Nashville, TN - January 14, 2017"  ->  "Chicago, IL - March 11, 2017"
#
# The text lives on Slide 1, in the second paragraph of the
# "Content Placeholder 2" shape (the first paragraph is "Brian Hansen").

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)

# Locate the shape that contains the old location/date text rather than
# hard-coding an index, in case shape ordering ever changes.
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasTextFrame -and $candidate.TextFrame.TextRange.Text -like "*Nashville, TN*") {
        $sh = $candidate
    }
}

$tr = $sh.TextFrame.TextRange

# The location/date text is the whole 2nd line (paragraph) of the box.
$line = $tr.Lines(2, 1)

# Rebuild the line chunk-by-chunk: set the line to the first chunk, then
# keep appending the remaining chunks onto the end of the text range. Each
# appended chunk becomes its own run, giving the same four-run split used
# in the updated deck ("Chicago, IL " / "<en-dash> " / "March 11, " / "2017").
$line.Text = "Chicago, IL "
$tr.Text = $tr.Text + [char]0x2013 + " "
$tr.Text = $tr.Text + "March 11, "
$tr.Text = $tr.Text + "2017"
